$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Add a new slide at the end of the deck using the "Title Only" layout
# (ppLayoutTitleOnly = 11) -- this is the layout whose title placeholder is
# typed "title" (not "ctrTitle"), matching the inserted slide.
# ---------------------------------------------------------------------------
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 11)

# ---------------------------------------------------------------------------
# Title placeholder: "Thanks for your attention!"
# ---------------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Left   = 0 / 12700
$title.Top    = 2460625 / 12700
$title.Width  = 12192000 / 12700
$title.Height = 104.37504237007875

$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Thanks for your attention!"
$titleRange.LanguageID = "en-GB"
$titleRange.Font.Size = 60
$titleRange.Font.Bold = $true
$titleRange.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------------
# "Link to code" textbox
# ---------------------------------------------------------------------------
$linkBox = $s.Shapes.AddTextbox(1, 3079750 / 12700, 6425168 / 12700, 9182100 / 12700, 400110 / 12700)
$linkBox.Name = "TextBox 8"

$linkFrame = $linkBox.TextFrame
$linkFrame.WordWrap = $true

$linkRange = $linkFrame.TextRange
$linkRange.Text = "Link to code: https://github.com/GitHub-User228/ContinuousMathematicalModelling"
$linkRange.LanguageID = "en-GB"
$linkRange.Font.Size = 20

$boldPart = $linkRange.Characters(1, 12)
$boldPart.Font.Bold = $true
